$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.847.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.316.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.86%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.21%  "

$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("E13").Value = "  +0.17%  "

$ws.Range("E14").Value = "  -1.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.676.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.317.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.788"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.771.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.09%  "

$ws.Range("E20").Value = "  +2.73%  "

$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("E22").Value = "  +0.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.96%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("E26").Value = "  +1.08%  "

$ws.Range("E27").Value = "  -2.49%  "

$ws.Range("E28").Value = "  +14.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("E30").Value = "  +1.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.33%  "

$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.69%  "

$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0697"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.95%  "

$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("E38").Value = "  +1.84%  "

$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("E40").Value = "  +1.32%  "

$ws.Range("E41").Value = "  -0.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.928.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.17%  "

$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("E46").Value = "  -2.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.92%  "

$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.544.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.05%  "

$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.12%  "
